# Se incorporan logs de Julio.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Julio")

# --- 1. Fill in the daily log values for July 9 - July 23 (rows 11-25), cols B..J ---
$data = @{}
$data[11] = @(81,7,1,3,1,0,0,1,0)
$data[12] = @(51,3,0,4,1,0,0,2,0)
$data[13] = @(73,5,1,2,0,1,0,0,0)
$data[14] = @(83,2,3,1,0,0,0,0,0)
$data[15] = @(87,2,4,1,0,0,0,0,0)
$data[16] = @(104,3,3,2,0,0,0,0,0)
$data[17] = @(78,0,4,0,1,1,0,0,0)
$data[18] = @(97,10,2,1,0,0,0,1,1)
$data[19] = @(75,4,2,2,0,1,0,3,2)
$data[20] = @(96,7,1,2,0,0,0,1,0)
$data[21] = @(57,3,3,1,0,0,0,1,0)
$data[22] = @(73,4,2,1,0,0,0,0,0)
$data[23] = @(79,3,1,2,1,0,0,0,0)
$data[24] = @(83,8,0,1,0,0,0,1,1)
$data[25] = @(46,7,4,2,0,1,0,3,0)

foreach ($row in 11..25) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# G19 picked up the underline style that the "Peaton Fallecido" column briefly carried that day
$ws.Range("G19").Font.Underline = $true

# --- 2. Turn on the table's Totals Row for "Tabla13" (Julio log) ---
$lo = $ws.ListObjects.Item("Tabla13")
$lo.ShowTotals = $true
$lo.ListColumns.Item("Fecha").TotalsRowLabel = "Total"
$lo.ListColumns.Item("Volcaduras").TotalsCalculation = -4157
$lo.ListColumns.Item("Peaton Atropellado").TotalsCalculation = -4157
$lo.ListColumns.Item("Motocilista Atropellado").TotalsCalculation = -4157
$lo.ListColumns.Item("Ciclista Atropellado").TotalsCalculation = -4157
$lo.ListColumns.Item("Peaton Fallecido").TotalsCalculation = -4157
$lo.ListColumns.Item("Ciclista Fallecido").TotalsCalculation = -4157
$lo.ListColumns.Item("Motociclista Fallecido").TotalsCalculation = -4157
$lo.ListColumns.Item("Automovilista Fallecido").TotalsCalculation = -4157

# --- 3. Restore the selection to where the user last left it ---
$ws.Range("A25:XFD25").Select()
